# Fix bugs in upload user bulk template:
# Rename the bulk-upload header columns on the "UserInfo" sheet so they
# match the field names expected by the importer (no spaces).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserInfo")

$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "PhoneNumber"
